$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# This edit corresponds to an (upstream, automated) content rebuild
# of the LOM3036 discipline sheet: the "Docentes responsaveis" /
# "Programa" / "Avaliacao" blocks (rows 13-24) got reshuffled, and
# the trailing two rows (25-26) were removed.
# ---------------------------------------------------------------

# --- Direct value assignments for cells whose final text is plain
#     (non-numeric / non-date-like) text -----------------------
$ws.Range("B10").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("C10").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("C15").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5840622 - Miguel Justino Ribeiro Barboza'
$ws.Range("C18").Value = '5840622 - Miguel Justino Ribeiro Barboza'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.'
$ws.Range("C19").Value = 'Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'A média do semestre será computada com base na relação:M=(P1+2P2)/3'
$ws.Range("C20").Value = 'A média do semestre será computada com base na relação:M=(P1+2P2)/3'
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2'
$ws.Range("C21").Value = 'A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2'
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B23").Value = 'LOM3013 -  Ciência dos Materiais  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOM3013 -  Ciência dos Materiais  (Requisito fraco)
'
$ws.Range("B24").Value = 'LOM3107 -  Mecânica dos Sólidos Deformáveis  (Requisito fraco)
'

# --- Clear cells that no longer hold content in the new layout -
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("A24").ClearContents()
$ws.Range("B25").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("B26").ClearContents()
$ws.Range("C26").ClearContents()

# --- B13/C13 must literally read "01/01/2018" as TEXT (shared
#     string), same as B8/C8 already do. Assigning that string
#     via .Value would get auto-parsed into a date serial, so
#     copy the already-text cells instead (values only, keeps
#     the existing column style s="2"/s="3" intact). -----------
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Drop the old trailing "Bibliografia" rows (25-26); the
#     sheet now ends at row 24 ------------------------------
$ws.Rows("25:26").Delete()

# --- Row heights for the reshuffled block -------------------
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(23).RowHeight = 30
$ws.Rows(24).RowHeight = 30

# --- Rows 17 and 22 revert to the sheet default row height ---
#     (no custom height any more) ---------------------------
$ws.Rows(17).AutoFit()
$ws.Rows(22).AutoFit()
